{"js": "// Append the new \"to dos\" content (date header + bullet list) to the\n// end of the document body, right after the final existing paragraph.\n// Office.js can't natively author a paragraph with zero inherited\n// formatting/list numbering in one call, so we build the exact target\n// OOXML for the 8 new paragraphs and insert it as a package fragment\n// after the current last paragraph -- this creates brand-new paragraphs\n// with only the formatting we specify (no bullet/numbering bleed-through\n// from the paragraph preceding them).\nconst lastParagraph = context.document.body.paragraphs.getLast();\nconst lastRange = lastParagraph.getRange(\"Whole\");\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>June 19</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n          <w:vertAlign w:val=\"superscript\"/>\n        </w:rPr>\n        <w:t>th</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: revise Stimuli and double-check items with James and Michael</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: agree on a date to record auditory stimuli with Roman</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: organize critical and filler items in folders</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: start writing the literature review and study design</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>To do perhaps: revise hypotheses (see Michael\u2019s mail)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>To do: start formulating the introductory part to the study, including fully acceptable, fully unacceptable, and neither acceptable nor unacceptable examples</w:t>\n      </w:r>\n    </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nlastRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Append the new \"to dos\" content (date header + bullet list) to the\n# end of the document, right after the last existing paragraph.\n# Word COM's Range can absorb arbitrary WordprocessingML via InsertXML\n# (the package/xmlData wrapper is the documented \"single file XML\"\n# format for Range.InsertXML / Range.WordOpenXML), so we build the\n# exact target OOXML for the 8 new paragraphs and insert it at a\n# collapsed range at the very end of the document. That creates clean\n# new paragraphs with only the formatting we specify -- no bullet or\n# numbering formatting bleeds through from the paragraph before them.\n$d = $word.ActiveDocument\n\n$endRange = $d.Range()\n$endRange.Collapse(0)\n\n$xml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>June 19</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n          <w:vertAlign w:val=\"superscript\"/>\n        </w:rPr>\n        <w:t>th</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: revise Stimuli and double-check items with James and Michael</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: agree on a date to record auditory stimuli with Roman</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: organize critical and filler items in folders</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Done: start writing the literature review and study design</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>To do perhaps: revise hypotheses (see Michael\u2019s mail)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Listenabsatz\"/>\n        <w:numPr>\n          <w:ilvl w:val=\"0\"/>\n          <w:numId w:val=\"1\"/>\n        </w:numPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>To do: start formulating the introductory part to the study, including fully acceptable, fully unacceptable, and neither acceptable nor unacceptable examples</w:t>\n      </w:r>\n    </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$endRange.InsertXML($xml) | Out-Null\n"}
